$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "cryptos" price/volume table (columns B-E, rows 2-51) with the
# latest scrape. A leading "'" forces plain-number-looking prices (e.g.
# "680.78") to remain stored as text, matching the original inline-string
# cells; values that already contain two dots (e.g. "69.246.49") are never
# auto-converted by Excel, so no prefix is required for those.
$ws.Range("D2").Value = "69.246.49"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "3.688.16"
$ws.Range("E3").Value = "  -2.93%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'680.78"
$ws.Range("E5").Value = "  -3.37%  "
$ws.Range("D6").Value = "'162.38"
$ws.Range("D7").Value = "3.686.19"
$ws.Range("E7").Value = "  -2.95%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.499"
$ws.Range("E9").Value = "  -4.03%  "
$ws.Range("E10").Value = "  -7.38%  "
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Value = "'0.450"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "'0.0000236"
$ws.Range("E13").Value = "  -6.78%  "
$ws.Range("D14").Value = "'33.52"
$ws.Range("E14").Value = "  -6.98%  "
$ws.Range("D15").Value = "4.312.59"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").Value = "3.692.27"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("D17").Value = "69.322.06"
$ws.Range("E17").Value = "  -1.99%  "
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").Value = "'16.34"
$ws.Range("E19").Value = "  -6.23%  "
$ws.Range("D20").Value = "'6.61"
$ws.Range("E20").Value = "  -7.30%  "
$ws.Range("D21").Value = "'482.43"
$ws.Range("E21").Value = "  -3.08%  "
$ws.Range("D22").Value = "'9.79"
$ws.Range("E22").Value = "  -7.68%  "
$ws.Range("D23").Value = "'0.666"
$ws.Range("E23").Value = "  -8.41%  "
$ws.Range("D24").Value = "'79.97"
$ws.Range("E24").Value = "  -5.47%  "
$ws.Range("D25").Value = "3.835.38"
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("E26").Value = "  -11.37%  "
$ws.Range("D27").Value = "'11.50"
$ws.Range("E27").Value = "  -4.89%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  -7.92%  "
$ws.Range("E30").Value = "  -10.51%  "
$ws.Range("D32").Value = "'2.11"
$ws.Range("E32").Value = "  -4.40%  "
$ws.Range("E33").Value = "  -6.63%  "
$ws.Range("D34").Value = "'27.04"
$ws.Range("E34").Value = "  -6.76%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.167"
$ws.Range("E35").Value = "  -3.21%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "3.652.45"
$ws.Range("E37").Value = "  -3.10%  "
$ws.Range("D38").Value = "'8.53"
$ws.Range("E38").Value = "  -6.20%  "
$ws.Range("D39").Value = "'6.01"
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("D40").Value = "'0.0941"
$ws.Range("E40").Value = "  -7.47%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  -5.98%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "'0.962"
$ws.Range("E44").Value = "  -7.20%  "
$ws.Range("D45").Value = "'158.66"
$ws.Range("E45").Value = "  -3.62%  "
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").Value = "'2.84"
$ws.Range("E47").Value = "  -13.00%  "
$ws.Range("D48").Value = "'0.000278"
$ws.Range("E48").Value = "  -14.70%  "
$ws.Range("D49").Value = "'388.68"
$ws.Range("E49").Value = "  -8.62%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'1.29"
$ws.Range("E50").Value = "  -4.68%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'8.09"
$ws.Range("E51").Value = "  -5.96%  "
